# Apply cryptos list update (prices/volumes refreshed by GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its numeric-looking values as plain text,
# matching the original inline-string cell contents (e.g. "28.471.46", "1.001").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.471.46"
$ws.Range("E2").Value = "  +2.64%  "

$ws.Range("D3").Value = "1.828.38"
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").Value = "315.48"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("D7").Value = "0.5147"
$ws.Range("E7").Value = "  -4.08%  "

$ws.Range("D8").Value = "0.3937"
$ws.Range("E8").Value = "  +3.17%  "

$ws.Range("D9").Value = "0.07735"
$ws.Range("E9").Value = "  +4.25%  "

$ws.Range("D10").Value = "41.96"
$ws.Range("E10").Value = "  +1.09%  "

$ws.Range("E11").Value = "  +2.58%  "

$ws.Range("D12").Value = "21.09"
$ws.Range("E12").Value = "  +3.91%  "

$ws.Range("D13").Value = "6.286"
$ws.Range("E13").Value = "  +1.45%  "

$ws.Range("D14").Value = "7.587"
$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16").Value = "1.824.54"
$ws.Range("E16").Value = "  +1.90%  "

$ws.Range("D17").Value = "93.58"
$ws.Range("E17").Value = "  +5.95%  "

$ws.Range("E18").Value = "  +2.39%  "

$ws.Range("D19").Value = "0.06619"
$ws.Range("E19").Value = "  +1.99%  "

$ws.Range("D20").Value = "17.72"
$ws.Range("E20").Value = "  +2.56%  "

$ws.Range("E21").Value = "  -0.18%  "

$ws.Range("D22").Value = "6.085"
$ws.Range("E22").Value = "  +2.68%  "

$ws.Range("D23").Value = "28.471.73"
$ws.Range("E23").Value = "  +2.43%  "

$ws.Range("D24").Value = "11.14"
$ws.Range("E24").Value = "  +0.42%  "

$ws.Range("E25").Value = "  +7.52%  "

$ws.Range("D26").Value = "157.23"
$ws.Range("E26").Value = "  +0.56%  "

$ws.Range("D27").Value = "2.447"
$ws.Range("E27").Value = "  +5.91%  "

$ws.Range("D28").Value = "20.63"
$ws.Range("E28").Value = "  +2.15%  "

$ws.Range("D29").Value = "2.036.82"
$ws.Range("E29").Value = "  +1.95%  "

$ws.Range("D30").Value = "124.93"
$ws.Range("E30").Value = "  +3.14%  "

$ws.Range("D31").Value = "1.132"
$ws.Range("E31").Value = "  +1.80%  "

$ws.Range("E32").Value = "  +0.50%  "

$ws.Range("D33").Value = "5.654"
$ws.Range("E33").Value = "  +2.82%  "

$ws.Range("D34").Value = "3.671"
$ws.Range("E34").Value = "  +0.63%  "

$ws.Range("D35").Value = "0.07154"
$ws.Range("E35").Value = "  +2.65%  "

$ws.Range("D36").Value = "0.2237"
$ws.Range("E36").Value = "  +1.93%  "

$ws.Range("D37").Value = "8.995"
$ws.Range("E37").Value = "  +6.64%  "

$ws.Range("E38").Value = "  +2.43%  "

$ws.Range("D39").Value = "5.154"
$ws.Range("E39").Value = "  +2.17%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6242"
$ws.Range("E40").Value = "  +2.44%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "11.27"
$ws.Range("E41").Value = "  -0.81%  "

$ws.Range("E42").Value = "  +2.64%  "

$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").Value = "1.397"
$ws.Range("E44").Value = "  -1.21%  "

$ws.Range("D45").Value = "13.43"
$ws.Range("E45").Value = "  +1.21%  "

$ws.Range("D46").Value = "0.5893"
$ws.Range("E46").Value = "  +3.69%  "

$ws.Range("D47").Value = "3.707"
$ws.Range("E47").Value = "  +0.83%  "

$ws.Range("D48").Value = "124.44"
$ws.Range("E48").Value = "  -0.03%  "

$ws.Range("D49").Value = "1.981"
$ws.Range("E49").Value = "  +4.13%  "

$ws.Range("E50").Value = "  +0.96%  "

$ws.Range("D51").Value = "0.06922"
$ws.Range("E51").Value = "  +1.91%  "
